$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion note text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$oldText = $ws1.Range("A1").Value()
$newText = $oldText.Replace(
    "1000 Bs = 13.79 = 56441.38 pesos",
    "1000 Bs = 13.66 = 55821.2 pesos"
).Replace(
    "56441.38 pesos = 13.77 = 978.78 Bs",
    "55821.2 pesos = 13.58 = 964.58 Bs"
)
$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update the rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 73.19799999999999
$ws2.Range("O10").Value = 4086
$ws2.Range("N12").Value = 4109.99
$ws2.Range("O12").Value = 71.02
